$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# StatDef sheet: monster stat table tweaks
# ---------------------------------------------------------------------------
$statDef = $wb.Worksheets.Item("StatDef")

# Class column (AC) corrections - several monsters had a placeholder "Normal"
# class that should have been their real race/class.
$statDef.Range("AC98").Value = "Dragon"
$statDef.Range("AC99").Value = "Dragon"
$statDef.Range("AC122").Value = "Strong"
$statDef.Range("AC135").Value = "Brute"
$statDef.Range("AC144").Value = "Golem,Elite"
$statDef.Range("AC235").Value = "Insect"
$statDef.Range("AC243").Value = "Brute"
$statDef.Range("AC245").Value = "Dragon"
$statDef.Range("AC246").Value = "Dragon"
$statDef.Range("AC251").Value = "Brute"
$statDef.Range("AC252").Value = "Demon"
$statDef.Range("AC253").Value = "Brute"

# Row 113 (Wild Rose) - stat multiplier tweaks
$statDef.Range("E113").Value = 90
$statDef.Range("H113").Value = 70
$statDef.Range("J113").Value = 130
$statDef.Range("L113").Value = 90
$statDef.Range("O113").Value = 80
$statDef.Range("R113").Value = 95

# Row 140 (Arclouse)
$statDef.Range("I140").Value = 110
$statDef.Range("L140").Value = 95

# Row 151 (Dokebi)
$statDef.Range("F151").Value = 130

# Row 152 (Am Mut)
$statDef.Range("F152").Value = 130

# Row 158 (Rideword)
$statDef.Range("E158").Value = 85
$statDef.Range("I158").Value = 130
$statDef.Range("K158").Value = 110
$statDef.Range("L158").Value = 105
$statDef.Range("O158").Value = 90
$statDef.Range("P158").Value = 90
$statDef.Range("Q158").Value = 102
$statDef.Range("R158").Value = 105

# Row 159 (Bathory)
$statDef.Range("E159").Value = 90
$statDef.Range("G159").Value = 110
$statDef.Range("H159").Value = 80
$statDef.Range("K159").Value = 110
$statDef.Range("O159").Value = 80
$statDef.Range("P159").Value = 110
$statDef.Range("R159").Value = 90

# Row 161 (Elder)
$statDef.Range("G161").Value = 130
$statDef.Range("R161").Value = 110

# Row 166 (Clock)
$statDef.Range("E166").Value = 120
$statDef.Range("L166").Value = 110
$statDef.Range("O166").Value = 90
$statDef.Range("P166").Value = 80
$statDef.Range("Q166").Value = 105
$statDef.Range("R166").Value = 110

# Restore view state: selection moved down to AC259 while keeping the frozen
# pane configuration intact.
$statDef.Activate()
$statDef.Range("AC259").Select() | Out-Null

# ---------------------------------------------------------------------------
# ClassDef sheet: insert a new "Dragon" class multiplier row
# ---------------------------------------------------------------------------
$classDef = $wb.Worksheets.Item("ClassDef")
$classDef.Activate()

$classDef.Rows.Item(15).Insert()

$classDef.Cells.Item(15, 1).Value = "Dragon"
$classDef.Cells.Item(15, 2).Value = 125
$classDef.Cells.Item(15, 3).Value = 120
$classDef.Cells.Item(15, 4).Value = 120
$classDef.Cells.Item(15, 5).Value = 110
$classDef.Cells.Item(15, 6).Value = 110
$classDef.Cells.Item(15, 7).Value = 90
$classDef.Cells.Item(15, 8).Value = 100
$classDef.Cells.Item(15, 9).Value = 100
$classDef.Cells.Item(15, 10).Value = 100
$classDef.Cells.Item(15, 11).Value = 100
$classDef.Cells.Item(15, 12).Value = 115
$classDef.Cells.Item(15, 13).Value = 105
$classDef.Cells.Item(15, 14).Value = 1

$classDef.Range("L15").Select() | Out-Null
